$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4
$ws.Range("I2").Value = 2.05
$ws.Range("L2").Value = 2.75
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.53
$ws.Range("U2").Value = 2.05
$ws.Range("V2").Value = 1.7
$ws.Range("W2").Value = 9
$ws.Range("X2").Value = 19
$ws.Range("AC2").Value = 7
$ws.Range("AK2").Value = 17
$ws.Range("AL2").Value = 19
$ws.Range("AM2").Value = 34
$ws.Range("BA2").Value = 67

# Row 3
$ws.Range("G3").Value = 2.3
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 3
$ws.Range("L3").Value = 4
$ws.Range("X3").Value = 10
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 21
$ws.Range("AC3").Value = 8
$ws.Range("AE3").Value = 15
$ws.Range("AG3").Value = 351
$ws.Range("AH3").Value = 9
$ws.Range("AO3").Value = 13
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 67

# Row 4
$ws.Range("G4").Value = 2.55
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 2.9
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 7

# Row 5
$ws.Range("G5").Value = 1.36
$ws.Range("L5").Value = 8.5
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.9
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.75
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("X5").Value = 6
$ws.Range("AD5").Value = 8.5
$ws.Range("AE5").Value = 23
$ws.Range("AH5").Value = 19
$ws.Range("AL5").Value = 67
$ws.Range("AM5").Value = 67
$ws.Range("AS5").Value = 151
$ws.Range("AT5").Value = 2.75
$ws.Range("AU5").Value = 10
$ws.Range("AZ5").Value = 201

# Row 10
$ws.Range("O10").Value = 1.14
$ws.Range("P10").Value = 5.5
